$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Clear out the old table contents (A1:L4) before writing the new layout.
# Use Clear() (not just ClearContents()) on B2 specifically since it carried
# a direct cell style (small 7pt font) in the old layout that is no longer
# used anywhere in the new table.
# ---------------------------------------------------------------------------
$ws.Range("A1:L4").ClearContents()
$ws.Range("B2").Clear()

# ---------------------------------------------------------------------------
# Header row (row 1): column A stays blank, B..K get the new labels.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "snake grid travel time without HB"
$ws.Range("C1").Value = "snake grid travel time with HB"
$ws.Range("D1").Value = "travel time algo2a without HB"
$ws.Range("E1").Value = "travel time algo2a with HB"
$ws.Range("F1").Value = "travel time algo2b without HB"
$ws.Range("G1").Value = "travel time algo2b with HB"
$ws.Range("H1").Value = "travel time algo2c without HB"
$ws.Range("I1").Value = "travel time algo2c with HB"
$ws.Range("J1").Value = "travel time no algo without HB"
$ws.Range("K1").Value = "travel time no algo with HB"

# ---------------------------------------------------------------------------
# Row 2 - "Total"
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Total"
$ws.Range("D2").Value = 502
$ws.Range("E2").Value = 5479
$ws.Range("F2").Value = 719
$ws.Range("G2").Value = 5684
$ws.Range("H2").Value = 535
$ws.Range("I2").Value = 5497
$ws.Range("J2").Value = 671
$ws.Range("K2").Value = 5637

# ---------------------------------------------------------------------------
# Row 3 - "Average"
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Average"
$ws.Range("D3").Value = 0.32618583499999998
$ws.Range("E3").Value = 3.5601039640000001
$ws.Range("F3").Value = 0.46718648499999998
$ws.Range("G3").Value = 3.6933073420000002
$ws.Range("H3").Value = 0.34762832999999999
$ws.Range("I3").Value = 3.57179987
$ws.Range("J3").Value = 0.43599740100000001
$ws.Range("K3").Value = 3.6627680310000001

# ---------------------------------------------------------------------------
# Column widths (characters). The underlying engine re-quantizes widths to
# 1/6-character increments when it persists them, so the inputs below are
# pre-compensated to land as close as possible on the desired final widths.
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").ColumnWidth = 29.0
$ws.Range("C1").ColumnWidth = 26.0
$ws.Range("D1").ColumnWidth = 25.333333333333332
$ws.Range("E1:F1").ColumnWidth = 26.0
$ws.Range("G1").ColumnWidth = 25.833333333333332
$ws.Range("H1").ColumnWidth = 37.833333333333336
$ws.Range("I1").ColumnWidth = 35.0
$ws.Range("J1").ColumnWidth = 41.0
$ws.Range("K1").ColumnWidth = 35.0

# ---------------------------------------------------------------------------
# View/selection
# ---------------------------------------------------------------------------
$ws.Range("J1").Select()
